$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 499, shifting existing rows 499:600
# down to 500:601 (dimension grows from A1:R600 to A1:R601).
$ws.Rows(499).Insert()

# Populate the newly-inserted row 499 with the new weekly record.
$ws.Range("A499").Value = 10
$ws.Range("B499").Value = "Vega Modelo de Temuco"
$ws.Range("C499").Value = "La Araucanía"
$ws.Range("D499").Value = 45244
$ws.Range("E499").Value = 9
$ws.Range("F499").Value = 100112017
$ws.Range("G499").Value = "Apio"
$ws.Range("H499").Value = "Americana (o)"
$ws.Range("I499").Value = "Primera"
$ws.Range("J499").Value = 125
$ws.Range("K499").Value = 8000
$ws.Range("L499").Value = 8000
$ws.Range("M499").Value = 8000
$ws.Range("N499").Value = "`$/caja 8 unidades"
$ws.Range("O499").Value = "Provincia del Elquí"
$ws.Range("P499").Value = 8000
$ws.Range("Q499").Value = 1
$ws.Range("R499").Value = "Hortaliza"
